$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A4").Value = 3
$ws.Range("B4").Value = 2
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = 3

$ws.Range("A6").Select()
